# Apply the commit's changes to the workbook.
#
# Summary of the edit (from the OOXML diff):
#  - Sheet "AISG" (sheet1): rows 2 and 4 drop two stale links
#    (gmail inbox link, peekingduck-trainer "tree/20221215" link) and ten
#    new reference links are appended after the existing ones, extending
#    the used range from A1:A14 to A1:A20.
#  - Sheet "Design Pattern" (sheet2): two now-orphaned "siim-isic" links
#    are removed from the shared-string pool as a side effect of the
#    sheet1 edits (its own cell content is unchanged).
#  - The active sheet/tab switches from "Design Pattern" back to "AISG",
#    with the AISG sheet's selection moved to E26.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AISG")
$ws2 = $wb.Worksheets.Item("Design Pattern")

# ---- Sheet "AISG": rewrite column A with the final link list ----
$links = @(
    "http://karpathy.github.io/2019/04/25/recipe/",
    "https://wordpress.deeplearning.ai/wp-content/uploads/2022/10/eBook-How-to-Build-a-Career-in-AI.pdf",
    "https://github.com/aisingapore/PeekingDuck-Familiarization-YOLOv6/blob/main/Instructions.md",
    "https://github.com/rwightman/pytorch-image-models/blob/main/inference.py",
    "https://cwiki.apache.org/confluence/display/MXNET/MXNet+-+Keras+Integration+Design",
    "https://eugeneyan.com/writing/design-patterns/",
    "https://github.com/msaroufim/ml-design-patterns",
    "https://refactoring.guru/design-patterns/observer/python/example#lang-features",
    "https://github.com/gao-hongnan/peekingduck-trainer/blob/255cda60aacf51e9c2cd7f36929330b6441bd338/src/dataset.py",
    "https://github.com/Lightning-AI/lightning/blob/master/examples/pl_loops/kfold.py",
    "https://torchmetrics.readthedocs.io/en/stable/classification/calibration_error.html",
    "https://torchmetrics.readthedocs.io/en/stable/pages/overview.html#metriccollection",
    "https://github.com/Lightning-AI/metrics/blob/master/src/torchmetrics/collections.py",
    "https://github.com/Lightning-AI/metrics/blob/96862e0d8175da57f39e573230f6878892882062/src/torchmetrics/metric.py#L44",
    "https://applyingml.com/resources/patterns/",
    "https://www.educative.io/courses/machine-learning-system-design?aid=5082902844932096&utm_source=google&utm_medium=paid&utm_campaign=machine-learning&utm_term=machine%20learning%20system%20design&utm_campaign=%5BTopics%5D+Machine+Learning&utm_source=adwords&utm_medium=ppc&hsa_acc=5451446008&hsa_cam=16394614703&hsa_grp=136977609347&hsa_ad=585425922123&hsa_src=g&hsa_tgt=kwd-302633670944&hsa_kw=machine%20learning%20system%20design&hsa_mt=b&hsa_net=adwords&hsa_ver=3&gclid=Cj0KCQiA4OybBhCzARIsAIcfn9lahOmeXFWLhJ7fqpMJFfE9ciRpM9-lnwvKlgrHw2Z_QOTLi7oZR5gaAhoAEALw_wcB",
    "https://medium.com/@upu1994/how-easy-is-making-custom-keras-callbacks-c771091602da",
    "https://pytorch-lightning.readthedocs.io/en/stable/common/trainer.html",
    "https://github.com/pytorch/vision/blob/main/references/classification/utils.py",
    "https://github.com/Lightning-AI/lightning/blob/master/src/pytorch_lightning/callbacks/model_checkpoint.py"
)

for ($i = 0; $i -lt $links.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 1).Value = $links[$i]
}

# ---- Switch the active tab back to "AISG" and set its selection ----
$ws1.Activate() | Out-Null
$ws1.Range("E26").Select() | Out-Null

# ---- Design Pattern sheet keeps its own selection, just no longer the active tab ----
$ws2.Range("B23").Select() | Out-Null
$ws1.Activate() | Out-Null
